# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 260, shifting the existing rows
# 260-269 down to 261-270 (dimension grows from A1:R269 to A1:R270).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 260; Excel's default insert
# behaviour copies formatting (incl. the date number format on column D)
# from the row above, so the existing data below is preserved intact.
$ws.Rows("260:260").Insert()

# Populate the newly inserted row 260 with the new reading.
$ws.Cells.Item(260, 1).Value  = 9
$ws.Cells.Item(260, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(260, 3).Value  = "Metropolitana"
$ws.Cells.Item(260, 4).Value  = 44706
$ws.Cells.Item(260, 5).Value  = 13
$ws.Cells.Item(260, 6).Value  = 300000001
$ws.Cells.Item(260, 7).Value  = "Rabanito"
$ws.Cells.Item(260, 8).Value  = "Sin especificar"
$ws.Cells.Item(260, 9).Value  = "Primera"
$ws.Cells.Item(260, 10).Value = 13000
$ws.Cells.Item(260, 11).Value = 2500
$ws.Cells.Item(260, 12).Value = 3000
$ws.Cells.Item(260, 13).Value = 2731
$ws.Cells.Item(260, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(260, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(260, 16).Value = 27
$ws.Cells.Item(260, 17).Value = 100
$ws.Cells.Item(260, 18).Value = "Hortaliza"
